$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "IsSourceOf"
$ws.Range("D1").Value = "IsDerivedFrom"

$ws.Range("B9").Select()
